$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion note text (cell A1) ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newNote = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 6.08 = 24703.45 pesos`n✅ 24703.45 pesos = 6.09 = 964.5 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newNote

# --- Sheet "tasas": update N10, O10, N12, O12 ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 164.39
$ws2.Range("O10").Value = 4061

$ws2.Range("N12").Value = 4059.6
$ws2.Range("O12").Value = 158.5
